# Case_0_161 vm_pu.xlsx update: slack/ext-grid voltage setpoint
# changed from 1.05 pu to 1.02 pu (column B), with the resulting
# per-bus voltage-magnitude results (columns C:F and I:N) recomputed
# for the 380 kV case. Columns A, G and H are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B_F = New-Object "object[,]" 24,5
$arr_B_F[0,0] = 1.02
$arr_B_F[0,1] = 1.029790378231177
$arr_B_F[0,2] = 1.032434797430825
$arr_B_F[0,3] = 1.037666915295026
$arr_B_F[0,4] = 1.044116600494014
$arr_B_F[1,0] = 1.02
$arr_B_F[1,1] = 1.031455340026064
$arr_B_F[1,2] = 1.03401324764155
$arr_B_F[1,3] = 1.039244068791146
$arr_B_F[1,4] = 1.046026200310985
$arr_B_F[2,0] = 1.02
$arr_B_F[2,1] = 1.032523344385266
$arr_B_F[2,2] = 1.035025638085954
$arr_B_F[2,3] = 1.040252915167518
$arr_B_F[2,4] = 1.047241185262533
$arr_B_F[3,0] = 1.02
$arr_B_F[3,1] = 1.032970131503664
$arr_B_F[3,2] = 1.035449130529505
$arr_B_F[3,3] = 1.040674274147368
$arr_B_F[3,4] = 1.047747074219561
$arr_B_F[4,0] = 1.02
$arr_B_F[4,1] = 1.033045020849787
$arr_B_F[4,2] = 1.035520113550262
$arr_B_F[4,3] = 1.040744861324232
$arr_B_F[4,4] = 1.047831729916726
$arr_B_F[5,0] = 1.02
$arr_B_F[5,1] = 1.0325293229846
$arr_B_F[5,2] = 1.035031305087657
$arr_B_F[5,3] = 1.040258556185227
$arr_B_F[5,4] = 1.047247964122297
$arr_B_F[6,0] = 1.02
$arr_B_F[6,1] = 1.030355013829105
$arr_B_F[6,2] = 1.032970119311596
$arr_B_F[6,3] = 1.038202361576885
$arr_B_F[6,4] = 1.044766268594874
$arr_B_F[7,0] = 1.02
$arr_B_F[7,1] = 1.026450511984433
$arr_B_F[7,2] = 1.029267837817265
$arr_B_F[7,3] = 1.034487958716023
$arr_B_F[7,4] = 1.040232445371715
$arr_B_F[8,0] = 1.02
$arr_B_F[8,1] = 1.023796046083306
$arr_B_F[8,2] = 1.026750248339813
$arr_B_F[8,3] = 1.031947920622538
$arr_B_F[8,4] = 1.037097979329409
$arr_B_F[9,0] = 1.02
$arr_B_F[9,1] = 1.022633876220492
$arr_B_F[9,2] = 1.025647865043605
$arr_B_F[9,3] = 1.030832333681935
$arr_B_F[9,4] = 1.035713255492515
$arr_B_F[10,0] = 1.02
$arr_B_F[10,1] = 1.022200226640991
$arr_B_F[10,2] = 1.02523650317517
$arr_B_F[10,3] = 1.030415537075149
$arr_B_F[10,4] = 1.035194695155057
$arr_B_F[11,0] = 1.02
$arr_B_F[11,1] = 1.022293335823582
$arr_B_F[11,2] = 1.025324827893197
$arr_B_F[11,3] = 1.030505051678372
$arr_B_F[11,4] = 1.035306120114996
$arr_B_F[12,0] = 1.02
$arr_B_F[12,1] = 1.022598071051835
$arr_B_F[12,2] = 1.025613900531133
$arr_B_F[12,3] = 1.030797930792768
$arr_B_F[12,4] = 1.035670477562689
$arr_B_F[13,0] = 1.02
$arr_B_F[13,1] = 1.022785566188589
$arr_B_F[13,2] = 1.025791756146839
$arr_B_F[13,3] = 1.030978061214921
$arr_B_F[13,4] = 1.035894409466145
$arr_B_F[14,0] = 1.02
$arr_B_F[14,1] = 1.023872902267084
$arr_B_F[14,2] = 1.026823147812767
$arr_B_F[14,3] = 1.03202162240292
$arr_B_F[14,4] = 1.03718929287952
$arr_B_F[15,0] = 1.02
$arr_B_F[15,1] = 1.024551506634015
$arr_B_F[15,2] = 1.027466800181903
$arr_B_F[15,3] = 1.032671971030577
$arr_B_F[15,4] = 1.03799412425406
$arr_B_F[16,0] = 1.02
$arr_B_F[16,1] = 1.024946096566892
$arr_B_F[16,2] = 1.027841053032357
$arr_B_F[16,3] = 1.033049794110506
$arr_B_F[16,4] = 1.038460923315393
$arr_B_F[17,0] = 1.02
$arr_B_F[17,1] = 1.025080434627928
$arr_B_F[17,2] = 1.027968465051317
$arr_B_F[17,3] = 1.033178366769403
$arr_B_F[17,4] = 1.038619643486898
$arr_B_F[18,0] = 1.02
$arr_B_F[18,1] = 1.024478826230301
$arr_B_F[18,2] = 1.027397864625299
$arr_B_F[18,3] = 1.032602351789082
$arr_B_F[18,4] = 1.03790804768927
$arr_B_F[19,0] = 1.02
$arr_B_F[19,1] = 1.022508388856925
$arr_B_F[19,2] = 1.02552882834141
$arr_B_F[19,3] = 1.030711752388384
$arr_B_F[19,4] = 1.035563300348752
$arr_B_F[20,0] = 1.02
$arr_B_F[20,1] = 1.021258084521656
$arr_B_F[20,2] = 1.024342744856706
$arr_B_F[20,3] = 1.029509039268743
$arr_B_F[20,4] = 1.034064649180086
$arr_B_F[21,0] = 1.02
$arr_B_F[21,1] = 1.02192199340492
$arr_B_F[21,2] = 1.024972564044898
$arr_B_F[21,3] = 1.030147967791974
$arr_B_F[21,4] = 1.034861455945284
$arr_B_F[22,0] = 1.02
$arr_B_F[22,1] = 1.024511671169881
$arr_B_F[22,2] = 1.027429017275826
$arr_B_F[22,3] = 1.032633814403907
$arr_B_F[22,4] = 1.037946950153745
$arr_B_F[23,0] = 1.02
$arr_B_F[23,1] = 1.02746881286982
$arr_B_F[23,2] = 1.030233505617484
$arr_B_F[23,3] = 1.035459268750421
$arr_B_F[23,4] = 1.041423978521735
$ws.Range("B2:F25").Value = $arr_B_F

$arr_I_N = New-Object "object[,]" 24,6
$arr_I_N[0,0] = 1.025530627837826
$arr_I_N[0,1] = 1.034935204090576
$arr_I_N[0,2] = 1.035240037679677
$arr_I_N[0,3] = 1.040457146806198
$arr_I_N[0,4] = 1.046888555230884
$arr_I_N[0,5] = 1.036404930789112
$arr_I_N[1,0] = 1.025746490092146
$arr_I_N[1,1] = 1.036237693180467
$arr_I_N[1,2] = 1.036625201758676
$arr_I_N[1,3] = 1.041842126006698
$arr_I_N[1,4] = 1.048606462676094
$arr_I_N[1,5] = 1.037709269562908
$arr_I_N[2,0] = 1.025875705303043
$arr_I_N[2,1] = 1.037070919323903
$arr_I_N[2,2] = 1.037512069345145
$arr_I_N[2,3] = 1.04272617204335
$arr_I_N[2,4] = 1.049697050643813
$arr_I_N[2,5] = 1.03854367898304
$arr_I_N[3,0] = 1.025927532733855
$arr_I_N[3,1] = 1.037418945245009
$arr_I_N[3,2] = 1.037882681659978
$arr_I_N[3,3] = 1.043094955147755
$arr_I_N[3,4] = 1.05015055390555
$arr_I_N[3,5] = 1.038892199140873
$arr_I_N[4,0] = 1.025936088770316
$arr_I_N[4,1] = 1.03747724848126
$arr_I_N[4,2] = 1.037944779356876
$arr_I_N[4,3] = 1.043156708109249
$arr_I_N[4,4] = 1.05022640849053
$arr_I_N[4,5] = 1.038950585174407
$arr_I_N[5,0] = 1.025876407613015
$arr_I_N[5,1] = 1.037075578514612
$arr_I_N[5,2] = 1.037517030195707
$arr_I_N[5,3] = 1.042731110973677
$arr_I_N[5,4] = 1.049703129878111
$arr_I_N[5,5] = 1.038548344790334
$arr_I_N[6,0] = 1.025605752124009
$arr_I_N[6,1] = 1.03537738616054
$arr_I_N[6,2] = 1.035710130553485
$arr_I_N[6,3] = 1.040927739051362
$arr_I_N[6,4] = 1.047473512424248
$arr_I_N[6,5] = 1.036847740808334
$arr_I_N[7,0] = 1.025048185766045
$arr_I_N[7,1] = 1.032310269153126
$arr_I_N[7,2] = 1.032452534312113
$arr_I_N[7,3] = 1.037655476869034
$arr_I_N[7,4] = 1.043381277903752
$arr_I_N[7,5] = 1.033776268142966
$arr_I_N[8,0] = 1.024621456462771
$arr_I_N[8,1] = 1.03021327089459
$arr_I_N[8,2] = 1.03022922401464
$arr_I_N[8,3] = 1.035408059546497
$arr_I_N[8,4] = 1.040539659534501
$arr_I_N[8,5] = 1.031676291906374
$arr_I_N[9,0] = 1.024423423685298
$arr_I_N[9,1] = 1.029292367624328
$arr_I_N[9,2] = 1.029253774775063
$arr_I_N[9,3] = 1.034418700953928
$arr_I_N[9,4] = 1.039281422383079
$arr_I_N[9,5] = 1.030754080847839
$arr_I_N[10,0] = 1.024347854467616
$arr_I_N[10,1] = 1.028948322670649
$arr_I_N[10,2] = 1.028889490260004
$arr_I_N[10,3] = 1.034048721726436
$arr_I_N[10,4] = 1.038809801088181
$arr_I_N[10,5] = 1.030409547310863
$arr_I_N[11,0] = 1.024364155652759
$arr_I_N[11,1] = 1.029022211832539
$arr_I_N[11,2] = 1.028967719931567
$arr_I_N[11,3] = 1.034128196946043
$arr_I_N[11,4] = 1.038911159377027
$arr_I_N[11,5] = 1.030483541403838
$arr_I_N[12,0] = 1.024417218252293
$arr_I_N[12,1] = 1.029263969424773
$arr_I_N[12,2] = 1.029223703139586
$arr_I_N[12,3] = 1.034388169406463
$arr_I_N[12,4] = 1.039242525317561
$arr_I_N[12,5] = 1.030725642319584
$arr_I_N[13,0] = 1.024449644784135
$arr_I_N[13,1] = 1.029412660470606
$arr_I_N[13,2] = 1.029381161841254
$arr_I_N[13,3] = 1.034548015823745
$arr_I_N[13,4] = 1.039446124391754
$arr_I_N[13,5] = 1.030874544523756
$arr_I_N[14,0] = 1.02463431832802
$arr_I_N[14,1] = 1.030274112942586
$arr_I_N[14,2] = 1.030293689132628
$arr_I_N[14,3] = 1.035473374093617
$arr_I_N[14,4] = 1.040622572219335
$arr_I_N[14,5] = 1.031737220357057
$arr_I_N[15,0] = 1.024746596848006
$arr_I_N[15,1] = 1.030810998111268
$arr_I_N[15,2] = 1.030862651111752
$arr_I_N[15,3] = 1.036049450769462
$arr_I_N[15,4] = 1.04135302970903
$arr_I_N[15,5] = 1.03227486796426
$arr_I_N[16,0] = 1.024810809375942
$arr_I_N[16,1] = 1.031122913615913
$arr_I_N[16,2] = 1.031193290823122
$arr_I_N[16,3] = 1.036383905953662
$arr_I_N[16,4] = 1.041776416751084
$arr_I_N[16,5] = 1.032587226424699
$arr_I_N[17,0] = 1.024832488045569
$arr_I_N[17,1] = 1.031229059627724
$arr_I_N[17,2] = 1.031305823768526
$arr_I_N[17,3] = 1.036497683351187
$arr_I_N[17,4] = 1.04192032916737
$arr_I_N[17,5] = 1.032693523176021
$arr_I_N[18,0] = 1.02473468271303
$arr_I_N[18,1] = 1.030753524039338
$arr_I_N[18,2] = 1.03080173402487
$arr_I_N[18,3] = 1.035987804979542
$arr_I_N[18,4] = 1.041274935934574
$arr_I_N[18,5] = 1.032217312272556
$arr_I_N[19,0] = 1.024401648316765
$arr_I_N[19,1] = 1.029192832874342
$arr_I_N[19,2] = 1.029148376921624
$arr_I_N[19,3] = 1.034311683088908
$arr_I_N[19,4] = 1.039145064454505
$arr_I_N[19,5] = 1.030654404747092
$arr_I_N[20,0] = 1.024180610506645
$arr_I_N[20,1] = 1.028200082617784
$arr_I_N[20,2] = 1.028097486291747
$arr_I_N[20,3] = 1.033243419921438
$arr_I_N[20,4] = 1.037781263272084
$arr_I_N[20,5] = 1.0296602446713
$arr_I_N[21,0] = 1.024298897597907
$arr_I_N[21,1] = 1.028727461910341
$arr_I_N[21,2] = 1.028655675633844
$arr_I_N[21,3] = 1.033811111126709
$arr_I_N[21,4] = 1.038506605911302
$arr_I_N[21,5] = 1.030188372902942
$arr_I_N[22,0] = 1.024740070144708
$arr_I_N[22,1] = 1.030779497929494
$arr_I_N[22,2] = 1.030829263620918
$arr_I_N[22,3] = 1.036015664878042
$arr_I_N[22,4] = 1.041310231415687
$arr_I_N[22,5] = 1.032243323048615
$arr_I_N[23,0] = 1.025201963616317
$arr_I_N[23,1] = 1.033112246131498
$arr_I_N[23,2] = 1.033303632189103
$arr_I_N[23,3] = 1.038512869240653
$arr_I_N[23,4] = 1.044458937502993
$arr_I_N[23,5] = 1.034579384020636
$ws.Range("I2:N25").Value = $arr_I_N
